$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "45.447.12"

# Row 3
$ws.Range("D3").Value = "2.415.16"
$ws.Range("E3").Value = "  +6.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.51"
$ws.Range("E5").Value = "  -2.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.70"
$ws.Range("E6").Value = "  -6.01%  "

# Row 7
$ws.Range("E7").Value = "  +1.20%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -0.75%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.97"
$ws.Range("E10").Value = "  -1.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  +0.96%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.05"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13
$ws.Range("E13").Value = "  +2.09%  "

# Row 14
$ws.Range("D14").Value = "2.790.56"
$ws.Range("E14").Value = "  +7.37%  "

# Row 15
$ws.Range("D15").Value = "2.421.27"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.837"
$ws.Range("E16").Value = "  +5.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.07"
$ws.Range("E17").Value = "  +4.01%  "

# Row 18
$ws.Range("D18").Value = "45.366.00"
$ws.Range("E18").Value = "  -3.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.35"
$ws.Range("E19").Value = "  -3.05%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0941"
$ws.Range("E20").Value = "  +1.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  +6.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.07"
$ws.Range("E22").Value = "  +2.64%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.99"
$ws.Range("E23").Value = "  -2.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.79"
$ws.Range("E24").Value = "  -0.43%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.10%  "

# Row 26
$ws.Range("E26").Value = "  +4.35%  "

# Row 27
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.79"
$ws.Range("E27").Value = "  -7.63%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.14"
$ws.Range("E30").Value = "  +5.21%  "

# Row 31
$ws.Range("E31").Value = "  +15.77%  "

# Row 32
$ws.Range("E32").Value = "  -2.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.80"
$ws.Range("E33").Value = "  +2.41%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("E34").Value = "  +1.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0764"
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.98"
$ws.Range("E36").Value = "  +17.64%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("E37").Value = "  -1.06%  "

# Row 38
$ws.Range("E38").Value = "  +0.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.84"
$ws.Range("E39").Value = "  -7.68%  "

# Row 40
$ws.Range("E40").Value = "  -2.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0298"
$ws.Range("E41").Value = "  +0.28%  "

# Row 42
$ws.Range("D42").Value = "2.008.57"
$ws.Range("E42").Value = "  +12.98%  "

# Row 43
$ws.Range("E43").Value = "  +3.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.12"
$ws.Range("E45").Value = "  -2.56%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.77"
$ws.Range("E46").Value = "  -7.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.13"
$ws.Range("E47").Value = "  +25.45%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.51"
$ws.Range("E48").Value = "  +9.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.83"
$ws.Range("E49").Value = "  +7.88%  "

# Row 50
$ws.Range("D50").Value = "2.659.26"
$ws.Range("E50").Value = "  +7.08%  "

# Row 51
$ws.Range("E51").Value = "  -0.77%  "
